# Refresh the "Price" (D) and "Volume(1h)" (E) columns for every
# coin row (2-51) with the latest scrape - GitHub Actions cron run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Source data stores these as plain text (e.g. "1.005", "22.20")
    # even though they look numeric. Assigning .Value directly lets
    # Excel auto-convert them to numbers, so round-trip the literal
    # through a text formula + paste-as-values instead, which keeps
    # the cell a string without leaving any NumberFormat residue.
    $range.Formula = '="' + ($text -replace '"', '""') + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $range.Worksheet.Application.CutCopyMode = $false
}

$ws.Range("D2").Value = '26.676.47'
$ws.Range("E2").Value = '  -6.93%  '
$ws.Range("D3").Value = '1.696.97'
$ws.Range("E3").Value = '  -5.67%  '
Set-TextValue $ws.Range("D4") '1.005'
$ws.Range("E4").Value = '  +0.21%  '
Set-TextValue $ws.Range("D5") '219.93'
$ws.Range("E5").Value = '  -5.07%  '
Set-TextValue $ws.Range("D6") '0.5131'
$ws.Range("E6").Value = '  -12.84%  '
Set-TextValue $ws.Range("D7") '1.005'
$ws.Range("E7").Value = '  +0.13%  '
Set-TextValue $ws.Range("D8") '0.2648'
$ws.Range("E8").Value = '  -4.38%  '
Set-TextValue $ws.Range("D9") '22.20'
$ws.Range("E9").Value = '  -4.46%  '
Set-TextValue $ws.Range("D10") '0.06302'
$ws.Range("E10").Value = '  -7.43%  '
Set-TextValue $ws.Range("D11") '0.07358'
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("D12").Value = '1.702.86'
$ws.Range("E12").Value = '  -5.97%  '
Set-TextValue $ws.Range("D13") '4.532'
$ws.Range("E13").Value = '  -5.26%  '
Set-TextValue $ws.Range("D14") '0.5803'
$ws.Range("E14").Value = '  -6.71%  '
$ws.Range("D15").Value = '1.928.53'
$ws.Range("E15").Value = '  -5.61%  '
Set-TextValue $ws.Range("D16") '0.000008468'
$ws.Range("E16").Value = '  -7.34%  '
Set-TextValue $ws.Range("D17") '65.57'
$ws.Range("E17").Value = '  -13.35%  '
$ws.Range("D18").Value = '26.713.24'
$ws.Range("E18").Value = '  -6.77%  '
Set-TextValue $ws.Range("D19") '5.013'
$ws.Range("E19").Value = '  -8.46%  '
Set-TextValue $ws.Range("D20") '1.005'
$ws.Range("E20").Value = '  +0.07%  '
Set-TextValue $ws.Range("D21") '11.01'
$ws.Range("E21").Value = '  -4.42%  '
Set-TextValue $ws.Range("D22") '186.92'
$ws.Range("E22").Value = '  -11.42%  '
Set-TextValue $ws.Range("D23") '6.256'
$ws.Range("E23").Value = '  -8.37%  '
Set-TextValue $ws.Range("D24") '1.006'
$ws.Range("E24").Value = '  +0.16%  '
Set-TextValue $ws.Range("D25") '144.58'
$ws.Range("E25").Value = '  -5.99%  '
Set-TextValue $ws.Range("D26") '7.525'
$ws.Range("E26").Value = '  -4.98%  '
Set-TextValue $ws.Range("D27") '0.1159'
$ws.Range("E27").Value = '  -8.70%  '
$ws.Range("E28").Value = '  -4.38%  '
Set-TextValue $ws.Range("D29") '1.353'
$ws.Range("E29").Value = '  -5.42%  '
Set-TextValue $ws.Range("D30") '0.05654'
$ws.Range("E30").Value = '  -7.68%  '
Set-TextValue $ws.Range("D31") '1.341'
$ws.Range("E31").Value = '  -5.83%  '
Set-TextValue $ws.Range("D32") '3.513'
$ws.Range("E32").Value = '  -7.19%  '
$ws.Range("E33").Value = '  -8.05%  '
$ws.Range("E34").Value = '  -5.04%  '
Set-TextValue $ws.Range("D35") '1.023'
$ws.Range("E35").Value = '  -3.31%  '
$ws.Range("E36").Value = '  -6.41%  '
Set-TextValue $ws.Range("D37") '2.358'
$ws.Range("E37").Value = '  -5.57%  '
Set-TextValue $ws.Range("D38") '2.699'
$ws.Range("E38").Value = '  -0.57%  '
Set-TextValue $ws.Range("D39") '0.01616'
$ws.Range("E39").Value = '  -4.95%  '
$ws.Range("D40").Value = '1.099.62'
$ws.Range("E40").Value = '  -4.26%  '
Set-TextValue $ws.Range("D41") '0.8603'
$ws.Range("E41").Value = '  -3.16%  '
Set-TextValue $ws.Range("D42") '5.849'
$ws.Range("E42").Value = '  -10.70%  '
$ws.Range("E43").Value = '  -0.41%  '
Set-TextValue $ws.Range("D44") '100.07'
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").Value = '1.856.97'
$ws.Range("E45").Value = '  -4.85%  '
$ws.Range("E46").Value = '  -0.94%  '
Set-TextValue $ws.Range("D47") '56.77'
$ws.Range("E47").Value = '  -5.89%  '
Set-TextValue $ws.Range("D48") '8.165'
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("E49").Value = '  +0.24%  '
Set-TextValue $ws.Range("D50") '0.05241'
$ws.Range("E50").Value = '  -4.06%  '
Set-TextValue $ws.Range("D51") '0.4320'
$ws.Range("E51").Value = '  -3.56%  '
